# Insert a new weekly price record at row 135 (Berenjena / Femacal de La Calera
# data set). All subsequent rows (old 135-143) shift down to 136-144, and the
# sheet's used range grows from A1:R143 to A1:R144.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 135..143 down to 136..144, leaving row 135 blank for the new record.
$ws.Rows.Item(135).Insert()

# Populate the new row 135 with the new data point.
$ws.Cells.Item(135, 1).Value  = 3
$ws.Cells.Item(135, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(135, 3).Value  = "Coquimbo"
$ws.Cells.Item(135, 4).Value  = 44461
$ws.Cells.Item(135, 5).Value  = 5
$ws.Cells.Item(135, 6).Value  = 100112001
$ws.Cells.Item(135, 7).Value  = "Berenjena"
$ws.Cells.Item(135, 8).Value  = "Sin especificar"
$ws.Cells.Item(135, 9).Value  = "Primera"
$ws.Cells.Item(135, 10).Value = 230
$ws.Cells.Item(135, 11).Value = 9000
$ws.Cells.Item(135, 12).Value = 9500
$ws.Cells.Item(135, 13).Value = 9261
$ws.Cells.Item(135, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(135, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(135, 16).Value = 154
$ws.Cells.Item(135, 17).Value = 60
$ws.Cells.Item(135, 18).Value = "Hortaliza"
